$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff removes the stored date values from A48:A51 (rows whose other
# columns already have no data), while keeping the cells' existing style
# (numFmt date style "s=2") intact. ClearContents clears the value/formula
# of the cells without touching formatting.
$ws.Range("A48:A51").ClearContents()
